$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D, shifting existing quarterly data from D:K to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy formatting (number format / style) from column F into the two new columns D and E
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)

# Update the dimension / values for the newly inserted & shifted columns
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(7, 6).Value = 43281
$ws.Cells.Item(7, 7).Value = 43190
$ws.Cells.Item(7, 8).Value = 43100
$ws.Cells.Item(7, 9).Value = 43008
$ws.Cells.Item(7, 10).Value = 42916
$ws.Cells.Item(7, 11).Value = 42825
$ws.Cells.Item(7, 12).Value = 42735
$ws.Cells.Item(7, 13).Value = 42643
$ws.Cells.Item(8, 4).Value = 344300
$ws.Cells.Item(8, 5).Value = 295200
$ws.Cells.Item(8, 6).Value = 759000
$ws.Cells.Item(8, 7).Value = 519700
$ws.Cells.Item(8, 8).Value = 285300
$ws.Cells.Item(8, 9).Value = 259100
$ws.Cells.Item(8, 10).Value = 619500
$ws.Cells.Item(8, 11).Value = 480900
$ws.Cells.Item(8, 12).Value = 498600
$ws.Cells.Item(8, 13).Value = 465700
$ws.Cells.Item(9, 4).Value = 192800
$ws.Cells.Item(9, 5).Value = 330300
$ws.Cells.Item(9, 6).Value = 450300
$ws.Cells.Item(9, 7).Value = 355300
$ws.Cells.Item(9, 8).Value = 519200
$ws.Cells.Item(9, 9).Value = 460300
$ws.Cells.Item(9, 10).Value = 655600
$ws.Cells.Item(9, 11).Value = 676600
$ws.Cells.Item(9, 12).Value = 680100
$ws.Cells.Item(9, 13).Value = 669900
$ws.Cells.Item(10, 4).Value = 151500
$ws.Cells.Item(10, 5).Value = -35100
$ws.Cells.Item(10, 6).Value = 308700
$ws.Cells.Item(10, 7).Value = 164500
$ws.Cells.Item(10, 8).Value = -233900
$ws.Cells.Item(10, 9).Value = -201200
$ws.Cells.Item(10, 10).Value = -36100
$ws.Cells.Item(10, 11).Value = -195700
$ws.Cells.Item(10, 12).Value = -181500
$ws.Cells.Item(10, 13).Value = -204300
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "NA"
$ws.Cells.Item(12, 8).Value = "NA"
$ws.Cells.Item(12, 9).Value = "NA"
$ws.Cells.Item(12, 10).Value = "NA"
$ws.Cells.Item(12, 11).Value = "NA"
$ws.Cells.Item(12, 12).Value = "NA"
$ws.Cells.Item(12, 13).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
$ws.Cells.Item(13, 9).Value = 0
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0
$ws.Cells.Item(13, 12).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 4).Value = -300
$ws.Cells.Item(14, 5).Value = -6200
$ws.Cells.Item(14, 6).Value = -16400
$ws.Cells.Item(14, 7).Value = -2000
$ws.Cells.Item(14, 8).Value = -18100
$ws.Cells.Item(14, 9).Value = -19800
$ws.Cells.Item(14, 10).Value = 500
$ws.Cells.Item(14, 11).Value = 200
$ws.Cells.Item(14, 12).Value = -900
$ws.Cells.Item(14, 13).Value = 100
$ws.Cells.Item(15, 4).Value = "NA"
$ws.Cells.Item(15, 5).Value = 9700
$ws.Cells.Item(15, 6).Value = 34100
$ws.Cells.Item(15, 7).Value = 18700
$ws.Cells.Item(15, 8).Value = 17400
$ws.Cells.Item(15, 9).Value = 15300
$ws.Cells.Item(15, 10).Value = 36100
$ws.Cells.Item(15, 11).Value = 18700
$ws.Cells.Item(15, 12).Value = 17800
$ws.Cells.Item(15, 13).Value = 15900
$ws.Cells.Item(17, 4).Value = 579800
$ws.Cells.Item(17, 5).Value = 72600
$ws.Cells.Item(17, 6).Value = 114100
$ws.Cells.Item(17, 7).Value = 457000
$ws.Cells.Item(17, 8).Value = 22900
$ws.Cells.Item(17, 9).Value = 168700
$ws.Cells.Item(17, 10).Value = 438700
$ws.Cells.Item(17, 11).Value = 464400
$ws.Cells.Item(17, 12).Value = 409200
$ws.Cells.Item(17, 13).Value = 393200
$ws.Cells.Item(18, 4).Value = -235500
$ws.Cells.Item(18, 5).Value = 222700
$ws.Cells.Item(18, 6).Value = 644900
$ws.Cells.Item(18, 7).Value = 62700
$ws.Cells.Item(18, 8).Value = 262400
$ws.Cells.Item(18, 9).Value = 90500
$ws.Cells.Item(18, 10).Value = 180700
$ws.Cells.Item(18, 11).Value = 16500
$ws.Cells.Item(18, 12).Value = 89400
$ws.Cells.Item(18, 13).Value = 72400
$ws.Cells.Item(20, 4).Value = -123100
$ws.Cells.Item(20, 5).Value = 25600
$ws.Cells.Item(20, 6).Value = -269600
$ws.Cells.Item(20, 7).Value = -4200
$ws.Cells.Item(20, 8).Value = -28100
$ws.Cells.Item(20, 9).Value = -64900
$ws.Cells.Item(20, 10).Value = 48200
$ws.Cells.Item(20, 11).Value = 48700
$ws.Cells.Item(20, 12).Value = 13100
$ws.Cells.Item(20, 13).Value = 6000
$ws.Cells.Item(21, 4).Value = -315500
$ws.Cells.Item(21, 5).Value = 274500
$ws.Cells.Item(21, 6).Value = 461100
$ws.Cells.Item(21, 7).Value = 85800
$ws.Cells.Item(21, 8).Value = 280900
$ws.Cells.Item(21, 9).Value = 45300
$ws.Cells.Item(21, 10).Value = 306500
$ws.Cells.Item(21, 11).Value = 95600
$ws.Cells.Item(21, 12).Value = 135500
$ws.Cells.Item(21, 13).Value = 107500
$ws.Cells.Item(22, 4).Value = "NA"
$ws.Cells.Item(22, 5).Value = 74200
$ws.Cells.Item(22, 6).Value = 176000
$ws.Cells.Item(22, 7).Value = 39800
$ws.Cells.Item(22, 8).Value = 45200
$ws.Cells.Item(22, 9).Value = 43400
$ws.Cells.Item(22, 10).Value = 139700
$ws.Cells.Item(22, 11).Value = 39300
$ws.Cells.Item(22, 12).Value = 42700
$ws.Cells.Item(22, 13).Value = 44700
$ws.Cells.Item(23, 4).Value = -358600
$ws.Cells.Item(23, 5).Value = 174100
$ws.Cells.Item(23, 6).Value = 199400
$ws.Cells.Item(23, 7).Value = 18700
$ws.Cells.Item(23, 8).Value = 189100
$ws.Cells.Item(23, 9).Value = -17900
$ws.Cells.Item(23, 10).Value = 89300
$ws.Cells.Item(23, 11).Value = 26000
$ws.Cells.Item(23, 12).Value = 59800
$ws.Cells.Item(23, 13).Value = 33700
$ws.Cells.Item(24, 4).Value = -92700
$ws.Cells.Item(24, 5).Value = 49600
$ws.Cells.Item(24, 6).Value = -2800
$ws.Cells.Item(24, 7).Value = 7800
$ws.Cells.Item(24, 8).Value = -115400
$ws.Cells.Item(24, 9).Value = 38200
$ws.Cells.Item(24, 10).Value = 63500
$ws.Cells.Item(24, 11).Value = 2300
$ws.Cells.Item(24, 12).Value = 11400
$ws.Cells.Item(24, 13).Value = 15500
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 4).Value = -265900
$ws.Cells.Item(26, 5).Value = 124500
$ws.Cells.Item(26, 6).Value = 202200
$ws.Cells.Item(26, 7).Value = 11000
$ws.Cells.Item(26, 8).Value = 304500
$ws.Cells.Item(26, 9).Value = -56100
$ws.Cells.Item(26, 10).Value = 25900
$ws.Cells.Item(26, 11).Value = 23700
$ws.Cells.Item(26, 12).Value = 48400
$ws.Cells.Item(26, 13).Value = 18200
$ws.Cells.Item(27, 4).Value = -164600
$ws.Cells.Item(27, 5).Value = 27200
$ws.Cells.Item(27, 6).Value = 57900
$ws.Cells.Item(27, 7).Value = 11600
$ws.Cells.Item(27, 8).Value = 50700
$ws.Cells.Item(27, 9).Value = 143600
$ws.Cells.Item(27, 10).Value = -24400
$ws.Cells.Item(27, 11).Value = 30500
$ws.Cells.Item(27, 12).Value = -25900
$ws.Cells.Item(27, 13).Value = 14400
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0
$ws.Cells.Item(28, 9).Value = 0
$ws.Cells.Item(28, 10).Value = 0
$ws.Cells.Item(28, 11).Value = 0
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(29, 4).Value = 17700
$ws.Cells.Item(29, 5).Value = -1300
$ws.Cells.Item(29, 6).Value = 286300
$ws.Cells.Item(29, 7).Value = -500
$ws.Cells.Item(29, 8).Value = 16800
$ws.Cells.Item(29, 9).Value = 12800
$ws.Cells.Item(29, 10).Value = 93900
$ws.Cells.Item(29, 11).Value = -31900
$ws.Cells.Item(29, 12).Value = 121100
$ws.Cells.Item(29, 13).Value = -9200
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 9).Value = 0
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 0
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 4).Value = 123100
$ws.Cells.Item(32, 5).Value = -25600
$ws.Cells.Item(32, 6).Value = 269600
$ws.Cells.Item(32, 7).Value = 4200
$ws.Cells.Item(32, 8).Value = 28100
$ws.Cells.Item(32, 9).Value = 64900
$ws.Cells.Item(32, 10).Value = -48200
$ws.Cells.Item(32, 11).Value = -48700
$ws.Cells.Item(32, 12).Value = -13100
$ws.Cells.Item(32, 13).Value = -6000
$ws.Cells.Item(33, 4).Value = -146900
$ws.Cells.Item(33, 5).Value = 25900
$ws.Cells.Item(33, 6).Value = 344200
$ws.Cells.Item(33, 7).Value = 11200
$ws.Cells.Item(33, 8).Value = 67500
$ws.Cells.Item(33, 9).Value = 156400
$ws.Cells.Item(33, 10).Value = 69500
$ws.Cells.Item(33, 11).Value = -1300
$ws.Cells.Item(33, 12).Value = 95200
$ws.Cells.Item(33, 13).Value = 5200
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 4).Value = -146900
$ws.Cells.Item(35, 5).Value = 25900
$ws.Cells.Item(35, 6).Value = 344200
$ws.Cells.Item(35, 7).Value = 11200
$ws.Cells.Item(35, 8).Value = 67500
$ws.Cells.Item(35, 9).Value = 156400
$ws.Cells.Item(35, 10).Value = 69500
$ws.Cells.Item(35, 11).Value = -1300
$ws.Cells.Item(35, 12).Value = 95200
$ws.Cells.Item(35, 13).Value = 5200
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(38, 6).Value = 43281
$ws.Cells.Item(38, 7).Value = 43190
$ws.Cells.Item(38, 8).Value = 43100
$ws.Cells.Item(38, 9).Value = 43008
$ws.Cells.Item(38, 10).Value = 42916
$ws.Cells.Item(38, 11).Value = 42825
$ws.Cells.Item(38, 12).Value = 42735
$ws.Cells.Item(38, 13).Value = 42643
$ws.Cells.Item(41, 4).Value = 1220800
$ws.Cells.Item(41, 5).Value = 1623900
$ws.Cells.Item(41, 6).Value = 1091200
$ws.Cells.Item(41, 7).Value = 756100
$ws.Cells.Item(41, 8).Value = 676600
$ws.Cells.Item(41, 9).Value = 597200
$ws.Cells.Item(41, 10).Value = 570200
$ws.Cells.Item(41, 11).Value = 585900
$ws.Cells.Item(41, 12).Value = 620500
$ws.Cells.Item(41, 13).Value = 404500
$ws.Cells.Item(42, 4).Value = 945400
$ws.Cells.Item(42, 5).Value = 1057900
$ws.Cells.Item(42, 6).Value = 876000
$ws.Cells.Item(42, 7).Value = 494800
$ws.Cells.Item(42, 8).Value = 500400
$ws.Cells.Item(42, 9).Value = 424900
$ws.Cells.Item(42, 10).Value = 327800
$ws.Cells.Item(42, 11).Value = 297600
$ws.Cells.Item(42, 12).Value = 309700
$ws.Cells.Item(42, 13).Value = 313000
$ws.Cells.Item(43, 4).Value = 442300
$ws.Cells.Item(43, 5).Value = 496000
$ws.Cells.Item(43, 6).Value = 448800
$ws.Cells.Item(43, 7).Value = 410300
$ws.Cells.Item(43, 8).Value = 394000
$ws.Cells.Item(43, 9).Value = 348400
$ws.Cells.Item(43, 10).Value = 365100
$ws.Cells.Item(43, 11).Value = 343500
$ws.Cells.Item(43, 12).Value = 367400
$ws.Cells.Item(43, 13).Value = 346100
$ws.Cells.Item(44, 4).Value = 17400
$ws.Cells.Item(44, 5).Value = 20200
$ws.Cells.Item(44, 6).Value = 18400
$ws.Cells.Item(44, 7).Value = 98700
$ws.Cells.Item(44, 8).Value = 96000
$ws.Cells.Item(44, 9).Value = 81000
$ws.Cells.Item(44, 10).Value = 97700
$ws.Cells.Item(44, 11).Value = 81600
$ws.Cells.Item(44, 12).Value = 87700
$ws.Cells.Item(44, 13).Value = 71900
$ws.Cells.Item(45, 4).Value = 337700
$ws.Cells.Item(45, 5).Value = 440700
$ws.Cells.Item(45, 6).Value = 374300
$ws.Cells.Item(45, 7).Value = 172300
$ws.Cells.Item(45, 8).Value = 166100
$ws.Cells.Item(45, 9).Value = 205000
$ws.Cells.Item(45, 10).Value = 141700
$ws.Cells.Item(45, 11).Value = 143200
$ws.Cells.Item(45, 12).Value = 150000
$ws.Cells.Item(45, 13).Value = 399800
$ws.Cells.Item(46, 4).Value = 2963600
$ws.Cells.Item(46, 5).Value = 3638600
$ws.Cells.Item(46, 6).Value = 2808700
$ws.Cells.Item(46, 7).Value = 1932200
$ws.Cells.Item(46, 8).Value = 1833200
$ws.Cells.Item(46, 9).Value = 1656500
$ws.Cells.Item(46, 10).Value = 1502400
$ws.Cells.Item(46, 11).Value = 1451800
$ws.Cells.Item(46, 12).Value = 1535300
$ws.Cells.Item(46, 13).Value = 1535300
$ws.Cells.Item(47, 4).Value = 950400
$ws.Cells.Item(47, 5).Value = 1114400
$ws.Cells.Item(47, 6).Value = 1058300
$ws.Cells.Item(47, 7).Value = 368300
$ws.Cells.Item(47, 8).Value = 340000
$ws.Cells.Item(47, 9).Value = 293700
$ws.Cells.Item(47, 10).Value = 300700
$ws.Cells.Item(47, 11).Value = 293000
$ws.Cells.Item(47, 12).Value = 274500
$ws.Cells.Item(47, 13).Value = 261100
$ws.Cells.Item(48, 4).Value = 5138500
$ws.Cells.Item(48, 5).Value = 5623700
$ws.Cells.Item(48, 6).Value = 5665600
$ws.Cells.Item(48, 7).Value = 3494900
$ws.Cells.Item(48, 8).Value = 3245600
$ws.Cells.Item(48, 9).Value = 2984800
$ws.Cells.Item(48, 10).Value = 2914900
$ws.Cells.Item(48, 11).Value = 2035700
$ws.Cells.Item(48, 12).Value = 1999300
$ws.Cells.Item(48, 13).Value = 1950500
$ws.Cells.Item(49, 4).Value = 352700
$ws.Cells.Item(49, 5).Value = 399200
$ws.Cells.Item(49, 6).Value = 433400
$ws.Cells.Item(49, 7).Value = 316900
$ws.Cells.Item(49, 8).Value = 293800
$ws.Cells.Item(49, 9).Value = 276600
$ws.Cells.Item(49, 10).Value = 284200
$ws.Cells.Item(49, 11).Value = 288700
$ws.Cells.Item(49, 12).Value = 295700
$ws.Cells.Item(49, 13).Value = 298800
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(50, 6).Value = 0
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 0
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 11).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 13).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = 0
$ws.Cells.Item(52, 4).Value = 404500
$ws.Cells.Item(52, 5).Value = 438600
$ws.Cells.Item(52, 6).Value = 498000
$ws.Cells.Item(52, 7).Value = 305300
$ws.Cells.Item(52, 8).Value = 264200
$ws.Cells.Item(52, 9).Value = 279200
$ws.Cells.Item(52, 10).Value = 302600
$ws.Cells.Item(52, 11).Value = 290500
$ws.Cells.Item(52, 12).Value = 244100
$ws.Cells.Item(52, 13).Value = 214900
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 9).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 11).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 13).Value = 0
$ws.Cells.Item(54, 4).Value = 9809600
$ws.Cells.Item(54, 5).Value = 11214500
$ws.Cells.Item(54, 6).Value = 9903000
$ws.Cells.Item(54, 7).Value = 6417600
$ws.Cells.Item(54, 8).Value = 5976700
$ws.Cells.Item(54, 9).Value = 5490800
$ws.Cells.Item(54, 10).Value = 5304700
$ws.Cells.Item(54, 11).Value = 4359700
$ws.Cells.Item(54, 12).Value = 4348800
$ws.Cells.Item(54, 13).Value = 4260600
$ws.Cells.Item(57, 4).Value = 325500
$ws.Cells.Item(57, 5).Value = 383800
$ws.Cells.Item(57, 6).Value = 431000
$ws.Cells.Item(57, 7).Value = 565900
$ws.Cells.Item(57, 8).Value = 530400
$ws.Cells.Item(57, 9).Value = 272400
$ws.Cells.Item(57, 10).Value = 291900
$ws.Cells.Item(57, 11).Value = 287200
$ws.Cells.Item(57, 12).Value = 309100
$ws.Cells.Item(57, 13).Value = 281200
$ws.Cells.Item(58, 4).Value = 1006900
$ws.Cells.Item(58, 5).Value = 1224100
$ws.Cells.Item(58, 6).Value = 748200
$ws.Cells.Item(58, 7).Value = 467200
$ws.Cells.Item(58, 8).Value = 442900
$ws.Cells.Item(58, 9).Value = 433300
$ws.Cells.Item(58, 10).Value = 457100
$ws.Cells.Item(58, 11).Value = 480200
$ws.Cells.Item(58, 12).Value = 487200
$ws.Cells.Item(58, 13).Value = 500100
$ws.Cells.Item(59, 4).Value = 191800
$ws.Cells.Item(59, 5).Value = 247100
$ws.Cells.Item(59, 6).Value = 191600
$ws.Cells.Item(59, 7).Value = 146000
$ws.Cells.Item(59, 8).Value = 129000
$ws.Cells.Item(59, 9).Value = 304600
$ws.Cells.Item(59, 10).Value = 316200
$ws.Cells.Item(59, 11).Value = 301100
$ws.Cells.Item(59, 12).Value = 281300
$ws.Cells.Item(59, 13).Value = 536300
$ws.Cells.Item(60, 4).Value = 1524200
$ws.Cells.Item(60, 5).Value = 1855000
$ws.Cells.Item(60, 6).Value = 1370800
$ws.Cells.Item(60, 7).Value = 1179100
$ws.Cells.Item(60, 8).Value = 1102300
$ws.Cells.Item(60, 9).Value = 1010200
$ws.Cells.Item(60, 10).Value = 1065200
$ws.Cells.Item(60, 11).Value = 1068600
$ws.Cells.Item(60, 12).Value = 1077600
$ws.Cells.Item(60, 13).Value = 1317600
$ws.Cells.Item(61, 4).Value = 5291900
$ws.Cells.Item(61, 5).Value = 6050800
$ws.Cells.Item(61, 6).Value = 5294200
$ws.Cells.Item(61, 7).Value = 3237800
$ws.Cells.Item(61, 8).Value = 2938300
$ws.Cells.Item(61, 9).Value = 2760600
$ws.Cells.Item(61, 10).Value = 2511700
$ws.Cells.Item(61, 11).Value = 2459200
$ws.Cells.Item(61, 12).Value = 2436100
$ws.Cells.Item(61, 13).Value = 2270000
$ws.Cells.Item(62, 4).Value = 901500
$ws.Cells.Item(62, 5).Value = 946200
$ws.Cells.Item(62, 6).Value = 991600
$ws.Cells.Item(62, 7).Value = 646000
$ws.Cells.Item(62, 8).Value = 610000
$ws.Cells.Item(62, 9).Value = 643900
$ws.Cells.Item(62, 10).Value = 641900
$ws.Cells.Item(62, 11).Value = 338400
$ws.Cells.Item(62, 12).Value = 335700
$ws.Cells.Item(62, 13).Value = 319900
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 0
$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 0
$ws.Cells.Item(64, 7).Value = 0
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 13).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 0
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 13).Value = 0
$ws.Cells.Item(66, 4).Value = 8838700
$ws.Cells.Item(66, 5).Value = 10051100
$ws.Cells.Item(66, 6).Value = 8750100
$ws.Cells.Item(66, 7).Value = 5714500
$ws.Cells.Item(66, 8).Value = 5275100
$ws.Cells.Item(66, 9).Value = 4891900
$ws.Cells.Item(66, 10).Value = 4711400
$ws.Cells.Item(66, 11).Value = 4264900
$ws.Cells.Item(66, 12).Value = 4270200
$ws.Cells.Item(66, 13).Value = 4248500
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 0
$ws.Cells.Item(69, 7).Value = 0
$ws.Cells.Item(69, 8).Value = 0
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 0
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 0
$ws.Cells.Item(69, 13).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 0
$ws.Cells.Item(70, 7).Value = 0
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = 0
$ws.Cells.Item(72, 4).Value = 593000
$ws.Cells.Item(72, 5).Value = 1128700
$ws.Cells.Item(72, 6).Value = 775100
$ws.Cells.Item(72, 7).Value = 668400
$ws.Cells.Item(72, 8).Value = 666900
$ws.Cells.Item(72, 9).Value = 564200
$ws.Cells.Item(72, 10).Value = 558600
$ws.Cells.Item(72, 11).Value = 55300
$ws.Cells.Item(72, 12).Value = 39100
$ws.Cells.Item(72, 13).Value = -27400
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0
$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 13).Value = 0
$ws.Cells.Item(76, 4).Value = 970900
$ws.Cells.Item(76, 5).Value = 1163400
$ws.Cells.Item(76, 6).Value = 1152900
$ws.Cells.Item(76, 7).Value = 703100
$ws.Cells.Item(76, 8).Value = 701600
$ws.Cells.Item(76, 9).Value = 598900
$ws.Cells.Item(76, 10).Value = 593300
$ws.Cells.Item(76, 11).Value = 94900
$ws.Cells.Item(76, 12).Value = 78600
$ws.Cells.Item(76, 13).Value = 12200
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(80, 6).Value = 43281
$ws.Cells.Item(80, 7).Value = 43190
$ws.Cells.Item(80, 8).Value = 43100
$ws.Cells.Item(80, 9).Value = 43008
$ws.Cells.Item(80, 10).Value = 42916
$ws.Cells.Item(80, 11).Value = 42825
$ws.Cells.Item(80, 12).Value = 42735
$ws.Cells.Item(80, 13).Value = 42643
$ws.Cells.Item(81, 4).Value = -146900
$ws.Cells.Item(81, 5).Value = 25900
$ws.Cells.Item(81, 6).Value = 344200
$ws.Cells.Item(81, 7).Value = 11200
$ws.Cells.Item(81, 8).Value = 67500
$ws.Cells.Item(81, 9).Value = 156400
$ws.Cells.Item(81, 10).Value = 69500
$ws.Cells.Item(81, 11).Value = -1300
$ws.Cells.Item(81, 12).Value = 95200
$ws.Cells.Item(81, 13).Value = 5200
$ws.Cells.Item(83, 4).Value = 43100
$ws.Cells.Item(83, 5).Value = 26200
$ws.Cells.Item(83, 6).Value = 85700
$ws.Cells.Item(83, 7).Value = 27300
$ws.Cells.Item(83, 8).Value = 38100
$ws.Cells.Item(83, 9).Value = 28400
$ws.Cells.Item(83, 10).Value = 108200
$ws.Cells.Item(83, 11).Value = 38100
$ws.Cells.Item(83, 12).Value = 39100
$ws.Cells.Item(83, 13).Value = 36900
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(84, 6).Value = 0
$ws.Cells.Item(84, 7).Value = 0
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(87, 6).Value = 0
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 11).Value = 0
$ws.Cells.Item(87, 12).Value = 0
$ws.Cells.Item(87, 13).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(88, 6).Value = 0
$ws.Cells.Item(88, 7).Value = 0
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = 0
$ws.Cells.Item(89, 4).Value = 68600
$ws.Cells.Item(89, 5).Value = 74600
$ws.Cells.Item(89, 6).Value = 321000
$ws.Cells.Item(89, 7).Value = -57100
$ws.Cells.Item(89, 8).Value = 205700
$ws.Cells.Item(89, 9).Value = 60600
$ws.Cells.Item(89, 10).Value = 204200
$ws.Cells.Item(89, 11).Value = 33700
$ws.Cells.Item(89, 12).Value = 59900
$ws.Cells.Item(89, 13).Value = 67400
$ws.Cells.Item(91, 4).Value = -23600
$ws.Cells.Item(91, 5).Value = -11300
$ws.Cells.Item(91, 6).Value = -43100
$ws.Cells.Item(91, 7).Value = -24800
$ws.Cells.Item(91, 8).Value = -14400
$ws.Cells.Item(91, 9).Value = -14200
$ws.Cells.Item(91, 10).Value = -29800
$ws.Cells.Item(91, 11).Value = -63400
$ws.Cells.Item(91, 12).Value = -77700
$ws.Cells.Item(91, 13).Value = -71900
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 13).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = 0
$ws.Cells.Item(94, 4).Value = -20900
$ws.Cells.Item(94, 5).Value = 24700
$ws.Cells.Item(94, 6).Value = -265500
$ws.Cells.Item(94, 7).Value = 31400
$ws.Cells.Item(94, 8).Value = -99400
$ws.Cells.Item(94, 9).Value = -126800
$ws.Cells.Item(94, 10).Value = -47400
$ws.Cells.Item(94, 11).Value = -34000
$ws.Cells.Item(94, 12).Value = 88900
$ws.Cells.Item(94, 13).Value = -32500
$ws.Cells.Item(96, 4).Value = 3200
$ws.Cells.Item(96, 5).Value = -5000
$ws.Cells.Item(96, 6).Value = -31900
$ws.Cells.Item(96, 7).Value = -31900
$ws.Cells.Item(96, 8).Value = 3000
$ws.Cells.Item(96, 9).Value = -3000
$ws.Cells.Item(96, 10).Value = 0
$ws.Cells.Item(96, 11).Value = -19000
$ws.Cells.Item(96, 12).Value = 9500
$ws.Cells.Item(96, 13).Value = -9500
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 9).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 13).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 0
$ws.Cells.Item(99, 7).Value = 0
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = 0
$ws.Cells.Item(100, 4).Value = -212300
$ws.Cells.Item(100, 5).Value = 231200
$ws.Cells.Item(100, 6).Value = -88700
$ws.Cells.Item(100, 7).Value = -216000
$ws.Cells.Item(100, 8).Value = 240400
$ws.Cells.Item(100, 9).Value = 92000
$ws.Cells.Item(100, 10).Value = 35300
$ws.Cells.Item(100, 11).Value = -40000
$ws.Cells.Item(100, 12).Value = 50800
$ws.Cells.Item(100, 13).Value = 6200
$ws.Cells.Item(101, 4).Value = -473600
$ws.Cells.Item(101, 5).Value = 437300
$ws.Cells.Item(101, 6).Value = 319100
$ws.Cells.Item(101, 7).Value = 41700
$ws.Cells.Item(101, 8).Value = 12200
$ws.Cells.Item(101, 9).Value = 1200
$ws.Cells.Item(101, 10).Value = 60000
$ws.Cells.Item(101, 11).Value = 5600
$ws.Cells.Item(101, 12).Value = 16200
$ws.Cells.Item(101, 13).Value = 600
$ws.Cells.Item(102, 4).Value = -638300
$ws.Cells.Item(102, 5).Value = 767800
$ws.Cells.Item(102, 6).Value = 285900
$ws.Cells.Item(102, 7).Value = 48700
$ws.Cells.Item(102, 8).Value = 110200
$ws.Cells.Item(102, 9).Value = 27000
$ws.Cells.Item(102, 10).Value = 252100
$ws.Cells.Item(102, 11).Value = -34600
$ws.Cells.Item(102, 12).Value = 216000
$ws.Cells.Item(102, 13).Value = 41400
